$d = $word.ActiveDocument

function Get-ParaNumber($doc, $rng) {
    $preceding = $doc.Range(0, $rng.Start)
    return $preceding.Paragraphs.Count + 1
}

function Find-ParaNumber($doc, $searchText) {
    $rng = $doc.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    return Get-ParaNumber $doc $rng
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$pkgFooter = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphInnerXml($doc, $paraNumber, $innerXml) {
    $p = $doc.Paragraphs.Item($paraNumber)
    $r = $p.Range
    $inner = $doc.Range($r.Start, $r.End - 1)
    $inner.InsertXML($pkgHeader + $innerXml + $pkgFooter)
}

# ---------------------------------------------------------------------------
# 1. Strip the redundant direct C00000 color formatting from the three
#    Heading4 paragraphs (Ingredients / Preparation / Accompagnement):
#    the color now lives purely in the Heading4 style definition.
# ---------------------------------------------------------------------------
foreach ($heading in @("Ingrédients", "Préparation", "Accompagnement")) {
    $n = Find-ParaNumber $d $heading
    $p = $d.Paragraphs.Item($n)
    $p.Style = "Heading4"
}

# ---------------------------------------------------------------------------
# 2. Remove the empty spacer paragraphs (pPr/ind left=0, no text) that sit
#    just before the "Preparation" and "Accompagnement" headings.
# ---------------------------------------------------------------------------
$n = Find-ParaNumber $d "1 bouquet garni, ou un bouillon cube de légumes"
$spacer = $d.Paragraphs.Item($n + 1)
if ($spacer.Range.Text.Trim().Length -eq 0) {
    $spacer.Range.Delete()
}

$n = Find-ParaNumber $d "Servir parsemé de persil."
$spacer = $d.Paragraphs.Item($n + 1)
if ($spacer.Range.Text.Trim().Length -eq 0) {
    $spacer.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3. Fix the "aire" -> "faire" typo, and drop the now-pointless proofErr
#    spell-check markers around it.
# ---------------------------------------------------------------------------
$n = Find-ParaNumber $d "Verser le mélange"
$inner = '<w:r><w:t xml:space="preserve">Verser le mélange sur la viande, et mélanger sans </w:t></w:r><w:r><w:t>faire</w:t></w:r><w:r><w:t xml:space="preserve"> bouillir.</w:t></w:r>'
Set-ParagraphInnerXml $d $n $inner

# ---------------------------------------------------------------------------
# 4. Merge the "Servir parsemé de persil." runs into a single run and drop
#    the grammar-check proofErr markers around "parsemé".
# ---------------------------------------------------------------------------
$n = Find-ParaNumber $d "Servir"
$inner = '<w:r><w:t>Servir parsemé de persil.</w:t></w:r>'
Set-ParagraphInnerXml $d $n $inner

# ---------------------------------------------------------------------------
# 5. Update the Heading4 / Heading4Char style colors from the themed orange
#    (E36C0A / accent6) to the plain dark red C00000.
# ---------------------------------------------------------------------------
$d.Styles.Item("Heading4").Font.Color = 192
$d.Styles.Item("Heading4Char").Font.Color = 192
